# Add the new "Company_ID" column (primary/foreign key column) to the
# student data sheet, matching the commit:
#   "Added company_id primary key and foreign key relationship between
#    company data and student data, updated upload function"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the Company_ID column (becomes shared string index 26,
# and widens the sheet dimension/row spans to column E automatically).
$ws.Range("E1").Value = "Company_ID"

# Match the new column's display width (closest value reachable through the
# ColumnWidth COM property, which Excel stores in quantized character units).
$ws.Columns.Item(5).ColumnWidth = 20.6

# Update the selection to the new column's data-entry range (E2:E5).
$ws.Range("E2:E5").Select()
